$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.333.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.063.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.29%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.055.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("E13").Value = "  +2.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.554.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.339.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.060.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.674"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  +3.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  +4.48%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "470.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0398"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.075.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("E43").Value = "  +4.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.255"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.11%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("E51").Value = "  +2.23%  "
